$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-08-05 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-08-06 Tuesday", 2) | Out-Null
$d.Content.Find.Execute("393÷9=43, 6", $true, $false, $false, $false, $false, $true, 1, $false, "275÷7=39, 2", 2) | Out-Null
$d.Content.Find.Execute("790÷9=87, 7", $true, $false, $false, $false, $false, $true, 1, $false, "621÷3=207, 0", 2) | Out-Null
$d.Content.Find.Execute("194÷9=21, 5", $true, $false, $false, $false, $false, $true, 1, $false, "655÷5=131, 0", 2) | Out-Null
$d.Content.Find.Execute("692÷6=115, 2", $true, $false, $false, $false, $false, $true, 1, $false, "412÷8=51, 4", 2) | Out-Null
$d.Content.Find.Execute("431÷4=107, 3", $true, $false, $false, $false, $false, $true, 1, $false, "323÷8=40, 3", 2) | Out-Null
$d.Content.Find.Execute("564÷7=80, 4", $true, $false, $false, $false, $false, $true, 1, $false, "542÷2=271, 0", 2) | Out-Null
$d.Content.Find.Execute("352÷7=50, 2", $true, $false, $false, $false, $false, $true, 1, $false, "732÷8=91, 4", 2) | Out-Null
$d.Content.Find.Execute("861÷5=172, 1", $true, $false, $false, $false, $false, $true, 1, $false, "829÷5=165, 4", 2) | Out-Null
$d.Content.Find.Execute("812÷6=135, 2", $true, $false, $false, $false, $false, $true, 1, $false, "215÷7=30, 5", 2) | Out-Null
$d.Content.Find.Execute("216÷9=24, 0", $true, $false, $false, $false, $false, $true, 1, $false, "723÷9=80, 3", 2) | Out-Null
$d.Content.Find.Execute("465÷3=155, 0", $true, $false, $false, $false, $false, $true, 1, $false, "793÷3=264, 1", 2) | Out-Null
$d.Content.Find.Execute("932÷7=133, 1", $true, $false, $false, $false, $false, $true, 1, $false, "602÷5=120, 2", 2) | Out-Null
$d.Content.Find.Execute("985÷2=492, 1", $true, $false, $false, $false, $false, $true, 1, $false, "721÷8=90, 1", 2) | Out-Null
$d.Content.Find.Execute("645÷3=215, 0", $true, $false, $false, $false, $false, $true, 1, $false, "640÷4=160, 0", 2) | Out-Null
$d.Content.Find.Execute("470÷5=94, 0", $true, $false, $false, $false, $false, $true, 1, $false, "674÷9=74, 8", 2) | Out-Null
$d.Content.Find.Execute("783÷6=130, 3", $true, $false, $false, $false, $false, $true, 1, $false, "146÷9=16, 2", 2) | Out-Null
$d.Content.Find.Execute("203÷7=29, 0", $true, $false, $false, $false, $false, $true, 1, $false, "471÷9=52, 3", 2) | Out-Null
$d.Content.Find.Execute("214÷9=23, 7", $true, $false, $false, $false, $false, $true, 1, $false, "162÷5=32, 2", 2) | Out-Null
$d.Content.Find.Execute("253÷9=28, 1", $true, $false, $false, $false, $false, $true, 1, $false, "201÷4=50, 1", 2) | Out-Null
$d.Content.Find.Execute("622÷4=155, 2", $true, $false, $false, $false, $false, $true, 1, $false, "736÷3=245, 1", 2) | Out-Null
$d.Content.Find.Execute("817÷3=272, 1", $true, $false, $false, $false, $false, $true, 1, $false, "502÷7=71, 5", 2) | Out-Null
$d.Content.Find.Execute("675÷2=337, 1", $true, $false, $false, $false, $false, $true, 1, $false, "619÷7=88, 3", 2) | Out-Null
$d.Content.Find.Execute("531÷4=132, 3", $true, $false, $false, $false, $false, $true, 1, $false, "628÷3=209, 1", 2) | Out-Null
$d.Content.Find.Execute("622÷8=77, 6", $true, $false, $false, $false, $false, $true, 1, $false, "641÷2=320, 1", 2) | Out-Null
$d.Content.Find.Execute("253÷3=84, 1", $true, $false, $false, $false, $false, $true, 1, $false, "284÷4=71, 0", 2) | Out-Null
